$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# --- Rename sheet "SwateTemplateMetadata" -> "isa_template" ---
$ws.Name = "isa_template"

# --- Remove column B (GEO / DPBO / DPBO ER tags) for rows 12-14, shifting
#     columns C:G left into B:F ---
$src = $ws.Range("C12:G14")
$scratch = $ws.Range("J12:N14")
$src.Copy($scratch)
$dst = $ws.Range("B12:G14")
$dst.Clear()
$scratch.Copy($ws.Range("B12:F14"))
$scratch.Clear()

# --- Re-point the hyperlink that used to live on G13 so it now lives on F13
#     (its text already moved there as part of the column shift above) ---
$ws.Range("G13").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F13"), "http://purl.obolibrary.org/obo/NCIT_C153189") | Out-Null

# Adding the hyperlink re-applies formatting and can create a duplicate
# "Link" style entry; nudging a format property back to its own value makes
# the engine collapse it back onto the original shared style.
$f13 = $ws.Range("F13")
$f13.Font.Underline = $f13.Font.Underline

# --- Row 12 wraps to fewer lines now that the long text lives in the wider
#     column B, so its height shrinks ---
$ws.Rows.Item(12).RowHeight = 28.8

# --- Update the active selection to reflect where the user ended up ---
$ws.Activate()
$ws.Range("B13").Select()

Write-Output "edit complete"
